$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.00"
$ws.Range("D4").Value = "'5.419"
$ws.Range("D5").Value = "'0.05897"
$ws.Range("D6").Value = "'3.451"
$ws.Range("D7").Value = "'6.555"
$ws.Range("D8").Value = "'0.8113"
$ws.Range("D9").Value = "'0.9427"
$ws.Range("D10").Value = "'0.1418"
$ws.Range("D11").Value = "'0.07435"
$ws.Range("D12").Value = "'0.03262"
$ws.Range("D13").Value = "'0.03056"
$ws.Range("D15").Value = "'3.876"
$ws.Range("D16").Value = "'0.001591"
$ws.Range("D17").Value = "'0.04677"
$ws.Range("D18").Value = "'0.0005965"
$ws.Range("D19").Value = "'0.005906"
$ws.Range("D20").Value = "'0.001266"
$ws.Range("D22").Value = "'0.00009509"
$ws.Range("D23").Value = "'3.599"
$ws.Range("D24").Value = "'2.132"
$ws.Range("D40").Value = "'0.03945"
$ws.Range("D41").Value = "'0.006189"
$ws.Range("D43").Value = "'0.002542"
$ws.Range("D44").Value = "'0.009144"
$ws.Range("D45").Value = "'0.00005202"
$ws.Range("D48").Value = "'0.002287"
